# Natmi following Dr Hou advice
# Updates the LR-pair statistics on Sheet1 to reflect the recomputed
# ligand/receptor-expressing cell counts (1 -> 3) and the resulting
# re-derived expression / specificity / edge-weight values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new value
$updates = @{
    2  = @{ E=3; F=1; G=63.211268; H=189.633804; I=0.4922609885657722; J=0.4922609885657722; K=3; L=1; M=0.5177623333333333; N=1.553287; O=0.08698878192060831; P=0.08698878192060831; Q=32.72841361263866; R=294.555722513748; S=0.04282118378237101; T=0.04282118378237101 }
    3  = @{ E=3; F=1; G=63.211268; H=189.633804; I=0.4922609885657722; J=0.4922609885657722; K=3; L=1; M=3.098536666666666; N=9.29561; O=0.5205823464105641; P=0.5205823464105641; Q=195.8624316444933; R=1762.76188480044; S=0.2562623804739536; T=0.2562623804739536 }
    4  = @{ E=3; F=1; G=63.211268; H=189.633804; I=0.4922609885657722; J=0.4922609885657722; K=3; L=1; M=2.335759666666667; N=7.007279; O=0.3924288716688277; P=0.3924288716688277; Q=147.6463302732573; R=1328.816972459316; S=0.1931774243094476; T=0.1931774243094477 }
    5  = @{ E=3; F=1; G=43.30706799999999; H=129.921204; I=0.3372560111523963; J=0.3372560111523963; K=3; L=1; M=0.5177623333333333; N=1.553287; O=0.08698878192060831; P=0.08698878192060831; Q=22.42276857750533; R=201.804917197548; S=0.02933748960555005; T=0.02933748960555005 }
    6  = @{ E=3; F=1; G=43.30706799999999; H=129.921204; I=0.3372560111523963; J=0.3372560111523963; K=3; L=1; M=3.098536666666666; N=9.29561; O=0.5205823464105641; P=0.5205823464105641; Q=134.1885381238266; R=1207.69684311444; S=0.1755695256267819; T=0.1755695256267819 }
    7  = @{ E=3; F=1; G=43.30706799999999; H=129.921204; I=0.3372560111523963; J=0.3372560111523963; K=3; L=1; M=2.335759666666667; N=7.007279; O=0.3924288716688277; P=0.3924288716688277; Q=101.1549027159907; R=910.394124443916; S=0.1323489959200645; T=0.1323489959200645 }
    8  = @{ E=3; F=1; G=21.891734; H=65.675202; I=0.1704830002818315; J=0.1704830002818315; K=3; L=1; M=0.5177623333333333; N=1.553287; O=0.08698878192060831; P=0.08698878192060831; Q=11.33471527655267; R=102.012437488974; S=0.01483010853268725; T=0.01483010853268724 }
    9  = @{ E=3; F=1; G=21.891734; H=65.675202; I=0.1704830002818315; J=0.1704830002818315; K=3; L=1; M=3.098536666666666; N=9.29561; O=0.5205823464105641; P=0.5205823464105641; Q=67.83234049591333; R=610.49106446322; S=0.0887504403098287; T=0.08875044030982869 }
    10 = @{ E=3; F=1; G=21.891734; H=65.675202; I=0.1704830002818315; J=0.1704830002818315; K=3; L=1; M=2.335759666666667; N=7.007279; O=0.3924288716688277; P=0.3924288716688277; Q=51.13382931059533; R=460.204463795358; S=0.06690245143931557; T=0.06690245143931556 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
